$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.165.54"
$ws.Range("E2").Value = "  -3.01%  "
$ws.Range("D3").Value = "1.600.51"
$ws.Range("E3").Value = "  -2.80%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'1.001"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "'301.76"
$ws.Range("E6").Value = "  -2.23%  "
$ws.Range("E7").Value = "  -3.00%  "
$ws.Range("D8").Value = "'0.3662"
$ws.Range("E8").Value = "  -4.23%  "
$ws.Range("D9").Value = "'47.60"
$ws.Range("E9").Value = "  -6.92%  "
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("E11").Value = "  -5.04%  "
$ws.Range("D12").Value = "'0.08090"
$ws.Range("E12").Value = "  -4.04%  "
$ws.Range("D13").Value = "'23.01"
$ws.Range("E13").Value = "  -3.53%  "
$ws.Range("D14").Value = "'6.637"
$ws.Range("E14").Value = "  -6.39%  "
$ws.Range("D15").Value = "'7.590"
$ws.Range("E15").Value = "  -2.41%  "
$ws.Range("D16").Value = "'0.00001267"
$ws.Range("E16").Value = "  -3.08%  "
$ws.Range("D17").Value = "1.597.52"
$ws.Range("E17").Value = "  -2.99%  "
$ws.Range("D18").Value = "'91.57"
$ws.Range("E18").Value = "  -3.15%  "
$ws.Range("D19").Value = "'0.06801"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("D20").Value = "'18.47"
$ws.Range("E20").Value = "  -6.01%  "
$ws.Range("D21").Value = "'6.602"
$ws.Range("E21").Value = "  -3.70%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'13.07"
$ws.Range("E23").Value = "  -3.58%  "
$ws.Range("D24").Value = "23.163.85"
$ws.Range("D25").Value = "'2.360"
$ws.Range("E25").Value = "  -4.83%  "
$ws.Range("D26").Value = "'2.911"
$ws.Range("E26").Value = "  -2.45%  "
$ws.Range("E27").Value = "  -3.98%  "
$ws.Range("D28").Value = "'151.20"
$ws.Range("E28").Value = "  -0.88%  "
$ws.Range("D29").Value = "'5.233"
$ws.Range("E29").Value = "  -3.48%  "
$ws.Range("D30").Value = "'131.96"
$ws.Range("E30").Value = "  -5.05%  "
$ws.Range("D31").Value = "'2.443"
$ws.Range("E31").Value = "  -1.66%  "
$ws.Range("D32").Value = "'7.105"
$ws.Range("E32").Value = "  -8.00%  "
$ws.Range("D33").Value = "1.774.11"
$ws.Range("E33").Value = "  -2.94%  "
$ws.Range("D34").Value = "'0.9835"
$ws.Range("E34").Value = "  -4.10%  "
$ws.Range("D35").Value = "'0.07730"
$ws.Range("E35").Value = "  -3.65%  "
$ws.Range("D36").Value = "'0.02786"
$ws.Range("E36").Value = "  -5.58%  "
$ws.Range("D37").Value = "'6.326"
$ws.Range("E37").Value = "  -5.46%  "
$ws.Range("D38").Value = "'0.2551"
$ws.Range("E38").Value = "  -4.84%  "
$ws.Range("D39").Value = "'0.08883"
$ws.Range("E39").Value = "  -2.30%  "
$ws.Range("D40").Value = "'10.10"
$ws.Range("E40").Value = "  -6.30%  "
$ws.Range("D41").Value = "'1.402"
$ws.Range("E41").Value = "  -1.14%  "
$ws.Range("D42").Value = "'0.7169"
$ws.Range("E42").Value = "  -4.61%  "
$ws.Range("E43").Value = "  -4.57%  "
$ws.Range("D44").Value = "'15.93"
$ws.Range("E44").Value = "  -2.03%  "
$ws.Range("D45").Value = "'0.6656"
$ws.Range("E45").Value = "  -3.57%  "
$ws.Range("D46").Value = "'2.315"
$ws.Range("E46").Value = "  -5.06%  "
$ws.Range("D47").Value = "'1.000"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("E48").Value = "  -2.47%  "
$ws.Range("D49").Value = "'0.07986"
$ws.Range("E49").Value = "  -3.52%  "
$ws.Range("D50").Value = "'131.70"
$ws.Range("E50").Value = "  -1.88%  "
$ws.Range("D51").Value = "'1.175"
$ws.Range("E51").Value = "  -3.66%  "
